# Fuel_selection sheet: add the missing "H2_Blend" fuel (hard-coded in the
# model equations, per commit message) and re-sort the fuel list A-Z so the
# new entry lands in its correct alphabetical spot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuel_selection")

# Row 32 is currently blank; "insert" it (no-op shift at the bottom of the
# used range) and populate it with the new fuel before sorting everything
# back into alphabetical order.
$ws.Rows.Item(32).Insert()
$ws.Range("A32").Value = "H2_Blend"
$ws.Range("B32").Value = 1

$sortRange = $ws.Range("A2:B32")
$sortKey = $ws.Range("A2")
$sortRange.Sort($sortKey)

# This sheet is now the one the author was last working in.
$ws.Activate()
$ws.Range("E12").Select()
